$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing values of row 19 (columns A through R) before shifting,
# so we can restore them (for the unaffected fields) into the newly inserted row.
$vals = @()
for ($c = 1; $c -le 18; $c++) {
    $vals += ,$ws.Cells.Item(19, $c).Value2()
}

# Insert a new blank row at position 19; this pushes the former rows 19-152
# down to become rows 20-153, growing the used range to A1:R153.
$ws.Rows("19:19").Insert()

# Re-populate the newly inserted row 19 with the data that used to live there,
# since most of its columns are unchanged.
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(19, $c).Value = $vals[$c - 1]
}

# Apply the updated Fecha (date serial) and Volumen values for this new record.
$ws.Cells.Item(19, 4).Value = 44490
$ws.Cells.Item(19, 10).Value = 6000
